$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 73469699
$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = '2018-10-02'
$ws.Range("AA2").Style = "Normal"
$ws.Range("AX2").Value = 'Jacob Rudhe'
$ws.Range("B2").Value = 89170
$ws.Range("E2").Value = 3215
$ws.Range("F2").Value = 'Rödgul trumpetsvamp'
$ws.Range("G2").Value = 'Craterellus lutescens'
$ws.Range("H2").Value = '(Fr.) Fr.'
$ws.Range("I2").Value = ""
$ws.Range("J2").ClearContents()
$ws.Range("Q2").Value = 556085.9585446554
$ws.Range("R2").Value = 6664742.03884617
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = '2018-10-02'
$ws.Range("Y2").Style = "Normal"

# Row 3
$ws.Range("A3").Value = 73469696
$ws.Range("B3").Value = 90008
$ws.Range("E3").Value = 6031
$ws.Range("F3").Value = 'Blomkålssvamp'
$ws.Range("G3").Value = 'Sparassis crispa'
$ws.Range("H3").Value = '(Wulfen:Fr.) Fr.'
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = '1'
$ws.Range("I3").Style = "Normal"
$ws.Range("J3").Value = 'fruktkroppar'
$ws.Range("Q3").Value = 556171.192173962
$ws.Range("R3").Value = 6664785.219463159

# Row 4
$ws.Range("A4").Value = 73469700
$ws.Range("B4").Value = 90645
$ws.Range("D4").Value = 'NT'
$ws.Range("E4").Value = 4361
$ws.Range("F4").Value = 'Orange taggsvamp'
$ws.Range("G4").Value = 'Hydnellum aurantiacum'
$ws.Range("H4").Value = '(Batsch:Fr.) P.Karst.'
$ws.Range("J4").Value = 'mycel'
$ws.Range("Q4").Value = 556085.9585446554
$ws.Range("R4").Value = 6664742.03884617

# Row 5
$ws.Range("A5").Value = 73469692
$ws.Range("B5").Value = 90665
$ws.Range("D5").Value = 'LC'
$ws.Range("E5").Value = 4366
$ws.Range("F5").Value = 'Skarp dropptaggsvamp'
$ws.Range("G5").Value = 'Hydnellum peckii'
$ws.Range("H5").Value = 'Banker'
$ws.Range("J5").ClearContents()
$ws.Range("Q5").Value = 555986.06257353
$ws.Range("R5").Value = 6664808.815289573

# Row 6
$ws.Range("A6").Value = 73469681
$ws.Range("B6").Value = 90671
$ws.Range("D6").Value = 'NT'
$ws.Range("E6").Value = 4368
$ws.Range("F6").Value = 'Dofttaggsvamp'
$ws.Range("G6").Value = 'Hydnellum suaveolens'
$ws.Range("H6").Value = '(Scop.:Fr.) P. Karst.'
$ws.Range("J6").Value = 'mycel'
$ws.Range("Q6").Value = 555993.2239594045
$ws.Range("R6").Value = 6664732.148454694

# Row 7
$ws.Range("A7").Value = 73469686
$ws.Range("B7").Value = 103265
$ws.Range("D7").Value = 'LC'
$ws.Range("E7").Value = 221144
$ws.Range("F7").Value = 'Grönpyrola'
$ws.Range("G7").Value = 'Pyrola chlorantha'
$ws.Range("H7").Value = 'Sw.'
$ws.Range("I7").Value = ""
$ws.Range("J7").ClearContents()
$ws.Range("K7").Value = 'fullt utvecklade blad'
$ws.Range("Q7").Value = 555578.2393848968
$ws.Range("R7").Value = 6664698.923133198

# Row 8
$ws.Range("A8").Value = 73469683
$ws.Range("B8").Value = 94121
$ws.Range("D8").Value = 'NT'
$ws.Range("E8").Value = 53
$ws.Range("F8").Value = 'Vedtrappmossa'
$ws.Range("G8").Value = 'Crossocalyx hellerianus'
$ws.Range("H8").Value = '(Nees ex Lindenb.) Meyl.'
$ws.Range("K8").ClearContents()
$ws.Range("Q8").Value = 555972.7712295325
$ws.Range("R8").Value = 6664697.935079473

# Row 9
$ws.Range("A9").Value = 73469685
$ws.Range("B9").Value = 103265
$ws.Range("D9").Value = 'LC'
$ws.Range("E9").Value = 221144
$ws.Range("F9").Value = 'Grönpyrola'
$ws.Range("G9").Value = 'Pyrola chlorantha'
$ws.Range("H9").Value = 'Sw.'
$ws.Range("K9").Value = 'fullt utvecklade blad'
$ws.Range("Q9").Value = 555772.9561257223
$ws.Range("R9").Value = 6664668.970280989

# Row 10
$ws.Range("A10").Value = 73469694
$ws.Range("B10").Value = 89392
$ws.Range("D10").Value = 'NT'
$ws.Range("E10").Value = 1202
$ws.Range("F10").Value = 'Ullticka'
$ws.Range("G10").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H10").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("I10").NumberFormat = "@"
$ws.Range("I10").Value = '1'
$ws.Range("I10").Style = "Normal"
$ws.Range("J10").Value = 'mycel'
$ws.Range("K10").ClearContents()
$ws.Range("Q10").Value = 556089.8366219404
$ws.Range("R10").Value = 6664847.789130096

# Row 11
$ws.Range("A11").Value = 73469703
$ws.Range("Q11").Value = 556151.1261185352
$ws.Range("R11").Value = 6664693.179299683

# Row 12
$ws.Range("A12").Value = 73469688
$ws.Range("B12").Value = 89356
$ws.Range("D12").Value = 'LC'
$ws.Range("E12").Value = 5447
$ws.Range("F12").Value = 'Vedticka'
$ws.Range("G12").Value = 'Fuscoporia viticola'
$ws.Range("H12").Value = '(Schwein.) Murrill'
$ws.Range("I12").Value = ""
$ws.Range("J12").ClearContents()
$ws.Range("Q12").Value = 555476.1431185424
$ws.Range("R12").Value = 6664349.871716912

# Row 13
$ws.Range("A13").Value = 73469702
$ws.Range("B13").Value = 90671
$ws.Range("D13").Value = 'NT'
$ws.Range("E13").Value = 4368
$ws.Range("F13").Value = 'Dofttaggsvamp'
$ws.Range("G13").Value = 'Hydnellum suaveolens'
$ws.Range("H13").Value = '(Scop.:Fr.) P. Karst.'
$ws.Range("I13").NumberFormat = "@"
$ws.Range("I13").Value = '1'
$ws.Range("I13").Style = "Normal"
$ws.Range("J13").Value = 'mycel'
$ws.Range("Q13").Value = 556125.2407155134
$ws.Range("R13").Value = 6664753.108494076

# Row 14
$ws.Range("A14").Value = 73628441
$ws.Range("AA14").NumberFormat = "@"
$ws.Range("AA14").Value = '2018-10-18'
$ws.Range("AA14").Style = "Normal"
$ws.Range("AX14").Value = 'Jacob Rudhe, Mårten Berglind'
$ws.Range("B14").Value = 90696
$ws.Range("E14").Value = 5448
$ws.Range("F14").Value = 'Svartvit taggsvamp'
$ws.Range("G14").Value = 'Phellodon connatus'
$ws.Range("H14").Value = '(Schultz) nom.prov'
$ws.Range("P14").Value = 'Rallvaråsen, Kallmora, Vstm'
$ws.Range("Q14").Value = 555669.8713961896
$ws.Range("R14").Value = 6664715.766085356
$ws.Range("Y14").NumberFormat = "@"
$ws.Range("Y14").Value = '2018-10-18'
$ws.Range("Y14").Style = "Normal"

# Row 15
$ws.Range("A15").Value = 73625715
$ws.Range("AF15").Value = ""
$ws.Range("B15").Value = 90665
$ws.Range("D15").Value = 'LC'
$ws.Range("E15").Value = 4366
$ws.Range("F15").Value = 'Skarp dropptaggsvamp'
$ws.Range("G15").Value = 'Hydnellum peckii'
$ws.Range("H15").Value = 'Banker'
$ws.Range("I15").Value = ""
$ws.Range("J15").Value = ""
$ws.Range("K15").Value = ""
$ws.Range("P15").Value = 'Stora Öfstjärnen, Vstm'
$ws.Range("Q15").Value = 555906.6440642555
$ws.Range("R15").Value = 6664777.194354862
$ws.Range("S15").Value = 25

# Row 16
$ws.Range("A16").Value = 73628443
$ws.Range("AF16").ClearContents()
$ws.Range("I16").NumberFormat = "@"
$ws.Range("I16").Value = '1'
$ws.Range("I16").Style = "Normal"
$ws.Range("J16").Value = 'mycel'
$ws.Range("K16").ClearContents()
$ws.Range("Q16").Value = 556124.0641905284
$ws.Range("R16").Value = 6664862.769563208
$ws.Range("S16").Value = 10
